# Update automatico via Actualizar 02-19-2021 12-45-36
# The "Fecha" (Date) column D holds an Excel serial date/time stamp that is
# refreshed on every automated run. Each block of 14 rows shares one
# timestamp; this commit shifts every block's timestamp up to the next one
# and introduces a fresh timestamp for the newest block (rows 2-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D15").Value  = 44246.53152002265
$ws.Range("D16:D29").Value = 44246.51022087963
$ws.Range("D30:D43").Value = 44246.48894626158
